$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.151151180267334
$ws.Range("B1").Value = 2.257282733917236
$ws.Range("C1").Value = 1.835969567298889
$ws.Range("D1").Value = 1.767743349075317
$ws.Range("E1").Value = 1.623548746109009
